# "Cash Flow Quantization Size" workbook: bump the CFQS sheet's
# Quantization Size ($) value in B2 from 100000 to 400000.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CFQS")
$ws.Range("B2").Value = 400000
